# Refresh computed Leve-profit figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# per-sheet, per-row, matching the latest scheduled market-data pull.
$wb = $excel.ActiveWorkbook

# A1-style column letter -> 1-based index, for the columns touched by this update (H..N).
$colIndex = @{ "H"=8; "I"=9; "J"=10; "K"=11; "L"=12; "M"=13; "N"=14 }

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{ Row=11; Set=@{"H"=1751.3077; "I"=1751.3077; "K"=1751.3077; "M"=-1611.3077}; Clear=@() },
    @{ Row=15; Set=@{"H"=3369090.2; "I"=3369090.2; "K"=10107270.6; "M"=-10107101.6}; Clear=@() },
    @{ Row=18; Set=@{"H"=6395.4; "J"=0; "L"=0}; Clear=@("N") },
    @{ Row=32; Set=@{"H"=0; "J"=0; "L"=0}; Clear=@("N") },
    @{ Row=40; Set=@{"H"=3962.2856; "I"=3770.5833; "K"=3770.5833; "M"=-3595.5833}; Clear=@() },
    @{ Row=87; Set=@{"H"=103326.664; "J"=103326.664; "L"=103326.664; "N"=-105822.664}; Clear=@() },
    @{ Row=90; Set=@{"H"=103326.664; "J"=103326.664; "L"=309979.992; "N"=-322459.992}; Clear=@() },
    @{ Row=98; Set=@{"H"=29412548; "I"=30303686; "K"=30303686; "M"=-30302188}; Clear=@() },
    @{ Row=107; Set=@{"H"=987; "I"=1024.6154; "K"=1024.6154; "M"=895.3846000000001}; Clear=@() },
    @{ Row=122; Set=@{"H"=29412548; "I"=30303686; "K"=90911058; "M"=-90908608}; Clear=@() },
    @{ Row=132; Set=@{"H"=782.16364; "I"=784.2593000000001; "J"=669; "K"=2352.7779; "L"=2007; "M"=177.2221; "N"=-7067}; Clear=@() },
    @{ Row=137; Set=@{"H"=635134.9; "I"=1253936.1; "J"=16333.625; "K"=3761808.3; "L"=49000.875; "M"=-3759258.3; "N"=-54100.875}; Clear=@() }
)
foreach ($u in $updates) {
    foreach ($col in $u.Set.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Set[$col]
    }
    foreach ($col in $u.Clear) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{ Row=32; Set=@{"H"=10351.16; "I"=3342.473; "J"=30298.96; "K"=3342.473; "L"=30298.96; "M"=-3055.473; "N"=-30872.96}; Clear=@() },
    @{ Row=61; Set=@{"H"=26790160; "I"=25003894; "J"=31255826; "K"=25003894; "L"=31255826; "M"=-25003682; "N"=-31256250}; Clear=@() },
    @{ Row=68; Set=@{"H"=0; "I"=0; "K"=0}; Clear=@("M") },
    @{ Row=71; Set=@{"H"=0; "I"=0; "K"=0}; Clear=@("M") },
    @{ Row=132; Set=@{"H"=15880330; "I"=25645558; "J"=11834.25; "K"=76936674; "L"=35502.75; "M"=-76934144; "N"=-40562.75}; Clear=@() },
    @{ Row=136; Set=@{"H"=26790160; "I"=25003894; "J"=31255826; "K"=75011682; "L"=93767478; "M"=-75009132; "N"=-93772578}; Clear=@() }
)
foreach ($u in $updates) {
    foreach ($col in $u.Set.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Set[$col]
    }
    foreach ($col in $u.Clear) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{ Row=64; Set=@{"H"=2544; "I"=4499.5; "K"=4499.5; "M"=-4274.5}; Clear=@() },
    @{ Row=67; Set=@{"H"=2544; "I"=4499.5; "K"=4499.5; "M"=-3719.5}; Clear=@() },
    @{ Row=86; Set=@{"H"=3283.5264; "J"=3500; "L"=3500; "N"=-5746}; Clear=@() },
    @{ Row=89; Set=@{"H"=3283.5264; "J"=3500; "L"=17500; "N"=-28732}; Clear=@() },
    @{ Row=94; Set=@{"H"=1494.6842; "I"=2113.625; "J"=1044.5454; "K"=2113.625; "L"=1044.5454; "M"=-1662.625; "N"=-1946.5454}; Clear=@() },
    @{ Row=134; Set=@{"H"=3864336.2; "I"=2181.6924; "K"=6545.0772; "M"=-4010.0772}; Clear=@() }
)
foreach ($u in $updates) {
    foreach ($col in $u.Set.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Set[$col]
    }
    foreach ($col in $u.Clear) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{ Row=31; Set=@{"H"=1020754.3; "I"=9173.789000000001; "J"=5825762; "K"=9173.789000000001; "L"=5825762; "M"=-8878.789000000001; "N"=-5826352}; Clear=@() },
    @{ Row=34; Set=@{"H"=1020754.3; "I"=9173.789000000001; "J"=5825762; "K"=9173.789000000001; "L"=5825762; "M"=-8971.789000000001; "N"=-5826166}; Clear=@() },
    @{ Row=86; Set=@{"H"=40134.965; "J"=81336.62; "L"=81336.62; "N"=-83582.62}; Clear=@() },
    @{ Row=89; Set=@{"H"=40134.965; "J"=81336.62; "L"=406683.1; "N"=-417915.1}; Clear=@() },
    @{ Row=132; Set=@{"H"=5364.9697; "I"=1685.7778; "J"=21921.334; "K"=5057.3334; "L"=65764.00199999999; "M"=-2527.3334; "N"=-70824.00199999999}; Clear=@() },
    @{ Row=141; Set=@{"H"=253182.1; "J"=272113.1; "L"=272113.1; "N"=-282473.1}; Clear=@() }
)
foreach ($u in $updates) {
    foreach ($col in $u.Set.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Set[$col]
    }
    foreach ($col in $u.Clear) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{ Row=58; Set=@{"H"=2162.25; "J"=1883; "L"=5649; "N"=-5905}; Clear=@() },
    @{ Row=131; Set=@{"H"=26471.5; "I"=17109.666; "J"=35833.332; "K"=51328.99800000001; "L"=107499.996; "M"=-46288.99800000001; "N"=-117579.996}; Clear=@() },
    @{ Row=141; Set=@{"H"=341308.78; "I"=434539.84; "K"=1303619.52; "M"=-1298439.52}; Clear=@() }
)
foreach ($u in $updates) {
    foreach ($col in $u.Set.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Set[$col]
    }
    foreach ($col in $u.Clear) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{ Row=3; Set=@{"H"=25000; "I"=0; "J"=25000; "K"=0; "L"=25000; "N"=-25224}; Clear=@("M") },
    @{ Row=15; Set=@{"H"=25000; "I"=0; "J"=25000; "K"=0; "L"=25000; "N"=-25340}; Clear=@("M") },
    @{ Row=22; Set=@{"H"=1141.1428; "I"=797.6; "K"=797.6; "M"=-502.6}; Clear=@() },
    @{ Row=23; Set=@{"H"=2869.5; "I"=2869.5; "K"=2869.5; "M"=-2639.5}; Clear=@() },
    @{ Row=27; Set=@{"H"=1141.1428; "I"=797.6; "K"=797.6; "M"=-690.6}; Clear=@() },
    @{ Row=40; Set=@{"H"=5322.7085; "I"=5879.769; "J"=4664.364; "K"=5879.769; "L"=4664.364; "M"=-5743.769; "N"=-4936.364}; Clear=@() },
    @{ Row=45; Set=@{"H"=34502.8; "I"=30020.5; "K"=30020.5; "M"=-29613.5}; Clear=@() },
    @{ Row=46; Set=@{"H"=3884.9285; "I"=2761.6667; "J"=4727.375; "K"=2761.6667; "L"=4727.375; "M"=-2573.6667; "N"=-5103.375}; Clear=@() },
    @{ Row=93; Set=@{"H"=34484300; "I"=55556964; "J"=1757.1818; "K"=55556964; "L"=1757.1818; "M"=-55555716; "N"=-4253.1818}; Clear=@() },
    @{ Row=132; Set=@{"H"=3783.1538; "I"=3783.1538; "J"=0; "K"=11349.4614; "L"=0; "M"=-8819.4614}; Clear=@("N") }
)
foreach ($u in $updates) {
    foreach ($col in $u.Set.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Set[$col]
    }
    foreach ($col in $u.Clear) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{ Row=45; Set=@{"H"=12309; "J"=12853.143; "L"=12853.143; "N"=-13835.143}; Clear=@() },
    @{ Row=81; Set=@{"H"=16199.4; "I"=7666.6665; "K"=15333.333; "M"=-14272.333}; Clear=@() },
    @{ Row=84; Set=@{"H"=16199.4; "I"=7666.6665; "K"=76666.66500000001; "M"=-71362.66500000001}; Clear=@() },
    @{ Row=113; Set=@{"H"=1168.25; "I"=114.5; "K"=343.5; "M"=1826.5}; Clear=@() },
    @{ Row=122; Set=@{"H"=2355.3076; "I"=2476.8333; "K"=7430.499899999999; "M"=-4980.499899999999}; Clear=@() }
)
foreach ($u in $updates) {
    foreach ($col in $u.Set.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Set[$col]
    }
    foreach ($col in $u.Clear) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}

